# 2017-01-31 update: energy.gov - chunk 7
# Rolls table_4_12_b.xlsx forward from the October 2016/2015 YTD edition to the
# November 2016/2015 YTD edition: updates the title/column headers and the
# petroleum-coke-cost figures (and their derived percentage changes) for the
# states/divisions whose data changed with the new month's release.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Title (row 1) ----
$ws.Range("A1").Value = "Table 4.12.B. Average Cost of Petroleum Coke Delivered for Electricity Generation by State, (Year-to-Date) November 2016 and 2015"

# ---- Column headers (row 4) ----
$ws.Range("B4").Value = "November 2016 YTD"
$ws.Range("E4").Value = "November 2016 YTD"
$ws.Range("G4").Value = "November 2016 YTD"

$ws.Range("C4").Value = "November 2015 YTD"
$ws.Range("F4").Value = "November 2015 YTD"
$ws.Range("H4").Value = "November 2015 YTD"

# ---- Row 16 (Vermont) ----
$ws.Range("E16").Value = 1.22
$ws.Range("F16").Value = 1.3

# ---- Row 21 (Pennsylvania) ----
$ws.Range("C21").Value = 1.67
$ws.Range("D21").Value = 0.024
$ws.Range("F21").Value = 1.67

# ---- Row 30 (Iowa) ----
$ws.Range("B30").Value = 1.53
$ws.Range("C30").Value = 2.17
$ws.Range("D30").Value = -0.29
$ws.Range("E30").Value = 1.53
$ws.Range("F30").Value = 2.17

# ---- Row 33 (West North Central) ----
$ws.Range("B33").Value = 1.53
$ws.Range("C33").Value = 2.17
$ws.Range("D33").Value = -0.29
$ws.Range("E33").Value = 1.53
$ws.Range("F33").Value = 2.17

# ---- Row 40 (Georgia) ----
$ws.Range("B40").Value = 1.55
$ws.Range("C40").Value = 1.69
$ws.Range("D40").Value = -0.083
$ws.Range("E40").Value = 1.55
$ws.Range("F40").Value = 1.69

# ---- Row 42 (South Atlantic) ----
$ws.Range("B42").Value = 1.55
$ws.Range("C42").Value = 1.69
$ws.Range("D42").Value = -0.083
$ws.Range("E42").Value = 1.55
$ws.Range("F42").Value = 1.69

# ---- Row 45 (Kentucky) ----
$ws.Range("B45").Value = 1.51
$ws.Range("C45").Value = 1.83
$ws.Range("D45").Value = -0.17
$ws.Range("E45").Value = 1.51
$ws.Range("F45").Value = 1.83

# ---- Row 47 (East South Central) ----
$ws.Range("B47").Value = 1.51
$ws.Range("C47").Value = 1.83
$ws.Range("D47").Value = -0.17
$ws.Range("E47").Value = 1.51
$ws.Range("F47").Value = 1.83

# ---- Row 66 (U.S. Total) ----
$ws.Range("B66").Value = 1.6
$ws.Range("C66").Value = 1.87
$ws.Range("D66").Value = -0.14
$ws.Range("E66").Value = 1.48
$ws.Range("F66").Value = 1.8
$ws.Range("H66").Value = 2.45
